$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# CT-96 (row 22) renamed to CT-100
$ws.Range("B22").Value = "CT-100"

# Insert a new row before row 14 (shifts rows 14-22 down to 15-23)
$ws.Rows.Item(14).Insert()

# Fill the new row 14 with the new test case (literal date/time serials, not text)
$ws.Range("B14").Value = "CT-44"
$ws.Cells.Item(14, 3).Value = 44540.448611111111
$ws.Cells.Item(14, 4).Value = 44540.490277777797
$ws.Range("C14").NumberFormat = "[h]:mm:ss;@"
$ws.Range("D14").NumberFormat = "[h]:mm:ss;@"
$ws.Range("E14").Value = "Failed"

# Rename CT-28 -> CL-29 (row 12, unaffected by the insert since it's above row 14)
$ws.Range("B12").Value = "CL-29"

# CT-23 (row 7) status Passed -> Failed
$ws.Range("E7").Value = "Failed"
$ws.Range("E7").Style = "Bad"

# CT-58 (was row 15, now row 16) status Passed -> Failed
$ws.Range("E16").Value = "Failed"
$ws.Range("E16").Style = "Bad"

# CT-71 (was row 19, now row 20) renamed to CT-73 (output/status unchanged)
$ws.Range("B20").Value = "CT-73"

# CT-74 (was row 21, now row 22) status Passed -> Failed
$ws.Range("E22").Value = "Failed"
$ws.Range("E22").Style = "Bad"

$ws.Range("B12").Select()
